$d = $word.ActiveDocument
Write-Host "Hello"
